$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin/Link (B/C) for the Decentraland <-> EnergySwap swap in rows 48/49
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"

# Update Price (D) column values - force Text storage to preserve exact formatting
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.954.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.885.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4585"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4067"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07963"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9900"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.897.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.901"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.046"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06553"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.009.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.411"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.207"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.115.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.098"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.396"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.002"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09313"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.603"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.403"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.272"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06045"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02213"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.257"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.178"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5768"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.256"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07467"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.256"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5444"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.890"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.12"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (E) column values
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  -2.35%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  -3.08%  "
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("E23").Value = "  -2.31%  "
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("E38").Value = "  -2.99%  "
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("E43").Value = "  -4.17%  "
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("E47").Value = "  +6.76%  "
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("E49").Value = "  -2.96%  "
$ws.Range("E50").Value = "  -4.28%  "
$ws.Range("E51").Value = "  -1.33%  "